$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.131.78"
$ws.Range("E2").Value = "  -0.36%  "
$ws.Range("D3").Value = "1.880.38"
$ws.Range("E3").Value = "  -1.49%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "313.54"
$ws.Range("E5").Value = "  -0.34%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").Value = "0.5082"
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("D8").Value = "0.3855"
$ws.Range("E8").Value = "  -1.94%  "
$ws.Range("D9").Value = "0.09119"
$ws.Range("E9").Value = "  -2.43%  "
$ws.Range("D10").Value = "1.122"
$ws.Range("E10").Value = "  -1.67%  "
$ws.Range("D11").Value = "41.54"
$ws.Range("E11").Value = "  -0.80%  "
$ws.Range("D12").Value = "6.350"
$ws.Range("E12").Value = "  -0.89%  "
$ws.Range("D13").Value = "20.77"
$ws.Range("E13").Value = "  -0.56%  "
$ws.Range("D14").Value = "1.883.92"
$ws.Range("E14").Value = "  -1.64%  "
$ws.Range("D15").Value = "7.205"
$ws.Range("E15").Value = "  -1.47%  "
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").Value = "0.00001112"
$ws.Range("E17").Value = "  -1.16%  "
$ws.Range("D18").Value = "91.22"
$ws.Range("E18").Value = "  -1.59%  "
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").Value = "18.21"
$ws.Range("E20").Value = "  +1.21%  "
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").Value = "6.116"
$ws.Range("E22").Value = "  -1.44%  "
$ws.Range("D23").Value = "28.172.21"
$ws.Range("E23").Value = "  -0.42%  "
$ws.Range("E24").Value = "  +0.23%  "
$ws.Range("D25").Value = "2.277"
$ws.Range("E25").Value = "  -1.71%  "
$ws.Range("D26").Value = "2.572"
$ws.Range("E26").Value = "  -0.92%  "
$ws.Range("D27").Value = "2.096.33"
$ws.Range("E27").Value = "  -1.81%  "
$ws.Range("D28").Value = "20.80"
$ws.Range("E28").Value = "  -1.24%  "
$ws.Range("D29").Value = "157.47"
$ws.Range("E29").Value = "  -0.39%  "
$ws.Range("D30").Value = "126.73"
$ws.Range("E30").Value = "  -0.47%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "0.1054"
$ws.Range("E31").Value = "  -1.74%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "1.064"
$ws.Range("E32").Value = "  -3.50%  "
$ws.Range("D33").Value = "5.616"
$ws.Range("E33").Value = "  -0.53%  "
$ws.Range("D34").Value = "3.603"
$ws.Range("E34").Value = "  -0.32%  "
$ws.Range("D35").Value = "9.689"
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("E36").Value = "  +1.81%  "
$ws.Range("D37").Value = "0.06579"
$ws.Range("E37").Value = "  -1.20%  "
$ws.Range("D38").Value = "0.2179"
$ws.Range("E38").Value = "  -0.37%  "
$ws.Range("D39").Value = "1.213"
$ws.Range("E39").Value = "  -2.76%  "
$ws.Range("D40").Value = "1.242"
$ws.Range("E40").Value = "  -2.49%  "
$ws.Range("D41").Value = "0.6416"
$ws.Range("E41").Value = "  -0.36%  "
$ws.Range("D42").Value = "11.56"
$ws.Range("E42").Value = "  +0.60%  "
$ws.Range("D43").Value = "4.928"
$ws.Range("E43").Value = "  -1.73%  "
$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "13.26"
$ws.Range("E45").Value = "  -0.15%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "0.6010"
$ws.Range("E46").Value = "  -0.16%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "3.676"
$ws.Range("E47").Value = "  -1.22%  "
$ws.Range("B48").Value = "WEMIXTOKEN"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").Value = "1.275"
$ws.Range("E48").Value = "  -0.16%  "
$ws.Range("D49").Value = "2.002"
$ws.Range("E49").Value = "  -0.87%  "
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").Value = "1.231"
$ws.Range("E50").Value = "  +3.72%  "
$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").Value = "121.33"
$ws.Range("E51").Value = "  -1.43%  "
